$d = $word.ActiveDocument
$sec = $d.Sections.First
$f2 = $sec.Footers.Item(1)
$ils = $f2.Range.InlineShapes.Item(1)
$shapeRange = $ils.Range
$xml = $shapeRange.WordOpenXML
$opts = [System.Text.RegularExpressions.RegexOptions]::Singleline
$pmatches = [regex]::Matches($xml, '<w:p\b.*?</w:p>', $opts)
$paraXml = $null
foreach ($pm in $pmatches) {
    if ($pm.Value.Contains('name="image2.png"')) {
        $paraXml = $pm.Value
    }
}
$paraXml = $paraXml.Replace('name="image2.png"', 'name="image1.png"')
$paraXml = [regex]::Replace($paraXml, 'r:embed="[^"]*"', 'r:embed="rId1"')

$pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$shapeRange.Text = ""
$shapeRange.InsertXML($pkg)
Write-Output "done"
